$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("pages_with_related_resources")
$ws2 = $wb.Worksheets.Item("pages_without_related_resources")
$ws3 = $wb.Worksheets.Item("pages_with_external_resources")

# --- Update the Spanish "vitamin D" path cells: node/116 -> node/146 ---
$ws1.Range("A4").Value = "espanol/node/146/2019/vitamina-d-complemento-cancer-prevencion"
$ws3.Range("A9").Value = "espanol/node/146/2019/vitamina-d-complemento-cancer-prevencion"

# --- Bold the header row on all three sheets ---
# sheet3's B1 already carries a left/top alignment style, so it ends up as a
# distinct (bold + alignment) style versus the other (bold only) header cells.
$ws3.Range("B1").Font.Bold = $true

$ws1.Range("A1:C1").Font.Bold = $true
$ws2.Range("A1:C1").Font.Bold = $true
$ws3.Range("A1").Font.Bold = $true
$ws3.Range("C1").Font.Bold = $true
$ws3.Range("D1").Font.Bold = $true

# --- Sheet view / selection updates ---
$ws2.Rows("1:1").Select()
$ws3.Rows("1:1").Select()

$ws1.Range("A18").Select()
$ws1.Activate()
